$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 757.15

$ws.Range("B5").Value = 529.374
$ws.Range("C5").Value = 545.8699999999999
$ws.Range("D5").Value = 545.51

$ws.Range("C6").Value = 338.7
$ws.Range("D6").Value = 306.6

$ws.Range("C7").Value = 432.4
$ws.Range("D7").Value = 426.2

$ws.Range("D8").Value = 443.2

$ws.Range("C9").Value = 747.788
$ws.Range("D9").Value = 747.1080000000001

$ws.Range("C10").Value = 1157.636
$ws.Range("D10").Value = 1157.036

$ws.Range("C11").Value = 1652.044
$ws.Range("D11").Value = 1651.444

$ws.Range("C12").Value = 1565.2
$ws.Range("D12").Value = 1564.6

$ws.Range("C13").Value = 649.7239999999999
$ws.Range("D13").Value = 649.124

$ws.Range("C14").Value = 257.732
$ws.Range("D14").Value = 257.052

$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0

$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0

$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0

$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0

$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0

$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0

$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0

$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
